# Update the Week 7 spreads tracker: shift rows 38-46 down into 39-47,
# and replace row 38 with the updated Marshall @ Old Dominion line.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row = 38; B = "Marshall";         C = "Old Dominion";    D = 5;               E = 1.199999999999999;  F = "Old Dominion -14.0"; G = "Old Dominion -14.5"; H = "Old Dominion -13.3"; I = -13.3;  J = -14.5 },
    @{ Row = 39; B = "Ole Miss";         C = "Washington State"; D = 3.7;             E = 1.100000000000001;  F = "Ole Miss -32.5";     G = "Ole Miss -32.5";     H = "Ole Miss -33.6";     I = 33.6;   J = 32.5 },
    @{ Row = 40; B = "Cincinnati";       C = "UCF";             D = 8.4;             E = 1.1;                 F = "Cincinnati -10.5";   G = "Cincinnati -11.0";   H = "Cincinnati -9.9";    I = 9.9;    J = 11 },
    @{ Row = 41; B = "Maryland";         C = "Nebraska";        D = 8.800000000000001; E = 0.9000000000000004; F = "Nebraska -5.5";      G = "Nebraska -6.5";      H = "Nebraska -5.6";      I = -5.6;   J = -6.5 },
    @{ Row = 42; B = "Notre Dame";       C = "NC State";        D = 7.7;             E = 0.8999999999999986; F = "Notre Dame -22.5";   G = "Notre Dame -21.5";   H = "Notre Dame -22.4";   I = 22.4;   J = 21.5 },
    @{ Row = 43; B = "Penn State";       C = "Northwestern";    D = 7.1;             E = 0.8999999999999986; F = "Penn State -22.5";   G = "Penn State -21.5";   H = "Penn State -20.6";   I = 20.6;   J = 21.5 },
    @{ Row = 44; B = "SMU";              C = "Stanford";        D = 5.7;             E = 0.8999999999999986; F = "SMU -18.0";          G = "SMU -19.5";          H = "SMU -18.6";          I = 18.6;   J = 19.5 },
    @{ Row = 45; B = "Texas Tech";       C = "Kansas";          D = 9.199999999999999; E = 0.8000000000000007; F = "Texas Tech -14.5";   G = "Texas Tech -13.5";   H = "Texas Tech -12.7";   I = 12.7;   J = 13.5 },
    @{ Row = 46; B = "Oklahoma State";   C = "Houston";         D = 4;               E = 0.8000000000000007; F = "Houston -13.5";      G = "Houston -14.0";      H = "Houston -14.8";      I = -14.8;  J = -14 },
    @{ Row = 47; B = "Georgia Southern"; C = "Southern Miss";   D = 5.4;             E = 0.7000000000000002; F = "Southern Miss -2.5"; G = "Southern Miss -3.5"; H = "Southern Miss -2.8"; I = -2.8;   J = -3.5 }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Range("B$row").Value = $r.B
    $ws.Range("C$row").Value = $r.C
    $ws.Range("D$row").Value = $r.D
    $ws.Range("E$row").Value = $r.E
    $ws.Range("F$row").Value = $r.F
    $ws.Range("G$row").Value = $r.G
    $ws.Range("H$row").Value = $r.H
    $ws.Range("I$row").Value = $r.I
    $ws.Range("J$row").Value = $r.J
}
